$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$meta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: add |4.0.1
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# Reference(...) value in Elements!K6 gets version suffix
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-document-reference|2.2.0-ballot)`n"

# Column K width grows to fit new (longer) content - closest value this
# host's pixel-quantised ColumnWidth setter can reach to the target 89.21875
$elements.Columns.Item(11).ColumnWidth = 88.25
